# New submission synced into the "JSS 3D" results sheet.
# Appends row 10 (Timestamp, Full Name, Admission No, AI Score) to the
# Google-Forms-style log table that lives in A1:D9 on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A10").Value = "2026-02-10 20:15:44"
$ws.Range("B10").Value = "Abubakar audu ali "
$ws.Range("C10").Value = "Number 7 "
$ws.Range("D10").Value = 10
